# "All terrains attached soma basic stuff done"
# Mark the "Basic terrain" (H) and "Complete terrain" (M) sub-columns as
# DONE for most rounds, set H16 to "není třeba", and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Cells that flip to "DONE" - copy the green "DONE" format (style used by
# B3, which is already filled in that way) onto each, then write the DONE
# value.
$doneCells = @("H3","M3","H4","M4","H5","M5","H6","M6","H9","H13","H14","H15","H19","H23","H24")

$ws.Range("B3").Copy()
foreach ($addr in $doneCells) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

foreach ($addr in $doneCells) {
    $ws.Range($addr).Value = "DONE"
}

# H16 becomes "není třeba" (gray), matching the style already used by G16.
$ws.Range("G16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H16").Value = "není třeba"

# C29 keeps its text ("POLOVINA") - no value change needed there; the shared
# string index shift in the file happens automatically because the two
# now-unused strings ("PŘIPRAVENO", "DOPLNIT POCITADLO LEVELU") are no
# longer referenced by any cell after the edits above.

# Update the view: drop the frozen/scrolled topLeftCell and move the
# selection to H5:H6.
$ws.Activate()
$ws.Range("H5:H6").Select()
